# Util function for delete temp files - rename sheet + append new employee rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "Cloud_EMS" to "EMS"
$ws.Name = "EMS"

$rows = @(
    @{ Row=2; EmpId=110; Name="David10"; Phone=9578821821; Age=30; Gender="Male"; Year=2010; Month=5; Day=21; Salary=50009 },
    @{ Row=3; EmpId=111; Name="David11"; Phone=9578821821; Age=30; Gender="Male"; Year=2010; Month=5; Day=22; Salary=50010 },
    @{ Row=4; EmpId=112; Name="David12"; Phone=9578821821; Age=30; Gender="Male"; Year=2010; Month=5; Day=23; Salary=50011 },
    @{ Row=5; EmpId=113; Name="David13"; Phone=9578821821; Age=30; Gender="Male"; Year=2010; Month=5; Day=24; Salary=50012 },
    @{ Row=6; EmpId=114; Name="David14"; Phone=9578821821; Age=30; Gender="Male"; Year=2010; Month=5; Day=25; Salary=50013 },
    @{ Row=7; EmpId=134; Name="David17"; Phone=9578821821; Age=30; Gender="Male"; Year=2010; Month=5; Day=28; Salary=50016 }
)

# Apply the built-in short-date number format (numFmtId 14) to the whole
# date column BEFORE any value is written there, so Excel reuses the
# built-in format id instead of minting a brand new custom numFmt. The
# format is set once on G2 and then propagated to G3:G7 via copy /
# paste-special so every cell shares the exact same style record (rather
# than each getting its own, duplicate, cellXfs entry).
$ws.Range("G2").NumberFormat = "mm-dd-yy"
$ws.Range("G2").Copy()
$ws.Range("G3:G7").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Range("A$rowNum").Value = $r.EmpId
    $ws.Range("B$rowNum").Value = $r.Name
    $ws.Range("C$rowNum").Value = $r.Phone
    $ws.Range("D$rowNum").Value = $r.Age
    $ws.Range("F$rowNum").Value = $r.Gender
    $ws.Range("G$rowNum").Value = (Get-Date -Year $r.Year -Month $r.Month -Day $r.Day -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
    $ws.Range("H$rowNum").Value = $r.Salary
}

Write-Host "Done. Sheet name: $($ws.Name)"
